$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for "ALL-VENT SYRUP 125ML" (row 5)
$ws.Rows(5).Delete()

# After that deletion, "NESTOGEN 1 MILK 400 GM" shifts up from row 25 to row 24
$ws.Rows(24).Delete()
